# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Thu Feb 15 09:13:45 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.938.78'
$ws.Range('E2').Value = '  +1.85%  '

$ws.Range('D3').Value = '2.781.00'
$ws.Range('E3').Value = '  +2.22%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '343.28'
$ws.Range('E5').Value = '  +3.83%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '115.51'
$ws.Range('E6').Value = '  +0.12%  '

$ws.Range('E7').Value = '  +3.23%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.576'
$ws.Range('E9').Value = '  +2.72%  '

$ws.Range('E10').Value = '  +4.27%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0852'
$ws.Range('E11').Value = '  +2.98%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.08'
$ws.Range('E12').Value = '  -1.23%  '

$ws.Range('E13').Value = '  +1.94%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.66'
$ws.Range('E14').Value = '  -0.01%  '

$ws.Range('D15').Value = '3.220.34'
$ws.Range('E15').Value = '  +2.47%  '

$ws.Range('D16').Value = '2.786.79'
$ws.Range('E16').Value = '  +3.24%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.882'
$ws.Range('E17').Value = '  +0.21%  '

$ws.Range('D18').Value = '51.851.54'
$ws.Range('E18').Value = '  +2.22%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.24'
$ws.Range('E19').Value = '  +9.07%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.04'
$ws.Range('E20').Value = '  +3.07%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.35'
$ws.Range('E21').Value = '  -3.61%  '

$ws.Range('D22').Value = '0.0₃0978'
$ws.Range('E22').Value = '  +1.76%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '269.96'
$ws.Range('E23').Value = '  -4.19%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.00'
$ws.Range('E24').Value = '  -0.21%  '

$ws.Range('E25').Value = '  +6.20%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.52'
$ws.Range('E26').Value = '  -1.03%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  -0.10%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.24'
$ws.Range('E28').Value = '  -1.05%  '

$ws.Range('E29').Value = '  +0.42%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.139'
$ws.Range('E30').Value = '  -1.30%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.65'
$ws.Range('E31').Value = '  -3.38%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '50.19'
$ws.Range('E32').Value = '  +0.42%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.71'
$ws.Range('E33').Value = '  +2.37%  '

$ws.Range('B34').Value = 'VeChain'
$ws.Range('C34').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0416'
$ws.Range('E34').Value = '  +18.16%  '

$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0820'
$ws.Range('E35').Value = '  -0.76%  '

$ws.Range('E36').Value = '  +0.18%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.96'
$ws.Range('E37').Value = '  -2.73%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.09'
$ws.Range('E38').Value = '  +0.41%  '

$ws.Range('E39').Value = '  -1.77%  '

$ws.Range('E40').Value = '  +0.39%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.68'
$ws.Range('E41').Value = '  +23.71%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '23.42'
$ws.Range('E42').Value = '  -1.53%  '

$ws.Range('E43').Value = '  +2.27%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '126.26'
$ws.Range('E44').Value = '  -2.02%  '

$ws.Range('E45').Value = '  -0.14%  '

$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '2.068.65'
$ws.Range('E46').Value = '  -1.93%  '

$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.33'
$ws.Range('E47').Value = '  -2.90%  '

$ws.Range('E48').Value = '  +1.22%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.53'
$ws.Range('E49').Value = '  +1.77%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.904'
$ws.Range('E50').Value = '  +13.34%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.90'
$ws.Range('E51').Value = '  -1.66%  '
